$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# Insert two new blank rows before the current row 25 so the existing
# row 25 data (2013-07-07 entry) moves down to row 27, leaving rows
# 25 and 26 free for new entries.
$ws.Rows("25:26").Insert()

# New row 25: 2013-07-03, 1h, "Revision manual"
$ws.Range("A25").Value = 41458
$ws.Range("B25").Value = 1
$ws.Range("D25").Value = "Revision manual"

# New row 26: 2013-07-06, 1h, "Revision manual"
$ws.Range("A26").Value = 41461
$ws.Range("B26").Value = 1
$ws.Range("D26").Value = "Revision manual"

# Copy the date style (style index 1, custom date format) used by column A
# down into the newly inserted cells so they match the rest of the column.
$ws.Range("A24").Copy()
$ws.Range("A25:A26").PasteSpecial(-4122) | Out-Null

# Append two further rows (28, 29) after the row that now holds the
# previously-existing 2013-07-07 entry (shifted to row 27).
$ws.Rows("28:29").Insert()

# New row 28: 2013-07-08, 2h, "Revision manual"
$ws.Range("A28").Value = 41463
$ws.Range("B28").Value = 2
$ws.Range("D28").Value = "Revision manual"

# New row 29: 2013-07-09, 1.5h, "Implementation tc14"
$ws.Range("A29").Value = 41464
$ws.Range("B29").Value = 1.5
$ws.Range("D29").Value = "Implementation tc14"

$ws.Range("A27").Copy()
$ws.Range("A28:A29").PasteSpecial(-4122) | Out-Null

$ws.Range("A29").Select()
